$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 303
$ws1.Range("F6").Value = 274
$ws1.Range("F10").Value = 53
$ws1.Range("F12").Value = 124
$ws1.Range("F13").Value = 2588
$ws1.Range("F19").Value = 548
$ws1.Range("F20").Value = 15
$ws1.Range("F21").Value = 627
$ws1.Range("F22").Value = 185
$ws1.Range("F26").Value = 56
$ws1.Range("F27").Value = 2222
$ws1.Range("F28").Value = 4780
$ws1.Range("F32").Value = 1242
$ws1.Range("F34").Value = 2156
$ws1.Range("D37").Value = "安源中大道17号 壹号公馆（萍乡）"
$ws1.Range("F38").Value = 54
$ws1.Range("F39").Value = 139
$ws1.Range("F41").Value = 444
$ws1.Range("F42").Value = 746
$ws1.Range("F43").Value = 13
$ws1.Range("F45").Value = 21
$ws1.Range("F46").Value = 439

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 49

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 303
$ws4.Range("F6").Value = 274
$ws4.Range("F10").Value = 53
$ws4.Range("F12").Value = 124
$ws4.Range("F13").Value = 2588
$ws4.Range("F17").Value = 49
$ws4.Range("F20").Value = 548
$ws4.Range("F21").Value = 15
$ws4.Range("F22").Value = 627
$ws4.Range("F23").Value = 185
$ws4.Range("F27").Value = 56
$ws4.Range("F28").Value = 2222
$ws4.Range("F29").Value = 4780
$ws4.Range("F33").Value = 1242
$ws4.Range("F35").Value = 2156
$ws4.Range("D38").Value = "安源中大道17号 壹号公馆（萍乡）"
$ws4.Range("F39").Value = 54
$ws4.Range("F40").Value = 139
$ws4.Range("F42").Value = 444
$ws4.Range("F43").Value = 746
$ws4.Range("F44").Value = 13
$ws4.Range("F46").Value = 21
$ws4.Range("F47").Value = 439
